$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Test: set row 7 data
$vals = @(1036736.6260914579, 1036583.2750369236, 1046252.4712205735, 1052006.7072076763)
for ($i = 0; $i -lt 4; $i++) {
    $col = 178 + $i
    $src = $ws.Cells.Item(7, 177)
    $dst = $ws.Cells.Item(7, $col)
    $src.Copy($dst)
    $dst.Value = $vals[$i]
}
$excel.CutCopyMode = $false
